$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.600.12'
$ws.Range("E2").Value = '  +4.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.792.85'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.93'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5356'
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3824'
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07531'
$ws.Range("E9").Value = '  +1.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  +2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.13'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.195'
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.413'
$ws.Range("E15").Value = '  +6.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.792.56'
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.47'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001066'
$ws.Range("E18").Value = '  +1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06440'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.31'
$ws.Range("E21").Value = '  +3.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.922'
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.613.17'
$ws.Range("E23").Value = '  +4.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.24'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.15'
$ws.Range("E26").Value = '  +3.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.58'
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.377'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.000.05'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.23'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.124'
$ws.Range("E31").Value = '  +4.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1017'
$ws.Range("E32").Value = '  -1.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.719'
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.652'
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2296'
$ws.Range("E35").Value = '  +11.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06573'
$ws.Range("E36").Value = '  +10.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02327'
$ws.Range("E37").Value = '  +2.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.097'
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.690'
$ws.Range("E39").Value = '  +5.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.50'
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6334'
$ws.Range("E41").Value = '  +3.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.207'
$ws.Range("E42").Value = '  +6.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.380'
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.55'
$ws.Range("E45").Value = '  +2.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5932'
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.671'
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.92'
$ws.Range("E48").Value = '  +3.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.983'
$ws.Range("E49").Value = '  +4.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.169'
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06925'
$ws.Range("E51").Value = '  +2.99%  '
